# Applies the "Updated test data and System overview document" edit to
# the "Scope & Phase" worksheet (sheet3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scope & Phase")
$ws.Activate()

# --- Highlight existing bullet items (yellow fill, same style already used for J6) ---
$ws.Range("D7").Interior.Color = 65535
$ws.Range("D14").Interior.Color = 65535

# --- Turn the old "Phase 1" second header (row 17) into the "Phase 2" header (row 18) ---
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Phase 2"
$ws.Range("B17").Clear()
$ws.Rows.Item(18).RowHeight = 14.25

# --- New bullet added under the Phase 2 section ---
$ws.Range("D28").Value = "Modernize application frontend"

# --- New "Extra" note row (red text) above the Phase 2 section header ---
$ws.Range("C16").Font.Color = 255
$ws.Range("C16").Value = "Extra "
$ws.Range("D16").Value = "Put something in the Home page????"

# --- Reorder the Phase 2 bullet list (Login form moves to the top) ---
$ws.Range("D20").Value = "Login form"
$ws.Range("D21").Value = "Generate order paper/pdf form"
$ws.Range("D22").Value = "Generate order fulfillment checklist"
$ws.Range("D23").Value = "Generate invoice paper "

# --- Update the saved selection/active cell ---
$ws.Range("K10").Select()

# --- Page setup touch-up ---
$ws.PageSetup.Orientation = 1
